$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 30 (B-SOC_DEV_ATH_M /
# "average training hours per male employee" row), shifting all the
# subsequent rows (old 30-44) down to (new 32-46).
$ws.Rows.Item(30).Insert()
$ws.Rows.Item(30).Insert()

# Populate the two new rows with the new Tag / Questions pairs.
# Values are written in this order so the shared-string table ends up with
# the new strings added in the same order as the authored workbook:
# "...turnover..." (B31), "...total number of employees..." (B30),
# "B-SOC-EMP-TNM" (A30), "B-SOC-EMP-TTN" (A31).
$ws.Range("B31").Value = "What is the total number and rate of employee turnover during the reporting period?"
$ws.Range("B30").Value = "What is the total number of employees at the end of the reporting period?"
$ws.Range("A30").Value = "B-SOC-EMP-TNM"
$ws.Range("A31").Value = "B-SOC-EMP-TTN"

# Match the row height (18) used by every other data row on the sheet.
$ws.Rows.Item(30).RowHeight = 18
$ws.Rows.Item(31).RowHeight = 18

# Update the active selection to match the post-edit view.
$ws.Range("A33").Select()
